# Apply the "cryptos list" GitHub Actions update: refreshed prices and
# 1h-volume percentages for every coin row, plus a swap of two ranking
# pairs (Kaspa/Binance-PegBSC-USD and OKB/InjectiveProtocol) that moved
# past each other in the source ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.967.72"
$ws.Range("E2").Value = "  -2.26%  "

$ws.Range("D3").Value = "3.482.81"
$ws.Range("E3").Value = "  -2.29%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'600.35"
$ws.Range("E5").Value = "  -3.13%  "

$ws.Range("D6").Value = "'147.92"
$ws.Range("E6").Value = "  -4.51%  "

$ws.Range("D7").Value = "3.480.44"
$ws.Range("E7").Value = "  -2.29%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  -2.43%  "

$ws.Range("D10").Value = "'0.142"
$ws.Range("E10").Value = "  -3.07%  "

$ws.Range("E11").Value = "  +3.50%  "

$ws.Range("D12").Value = "'0.423"
$ws.Range("E12").Value = "  -3.45%  "

$ws.Range("D13").Value = "'0.0000213"
$ws.Range("E13").Value = "  -4.08%  "

$ws.Range("D14").Value = "4.069.69"
$ws.Range("E14").Value = "  -2.35%  "

$ws.Range("D15").Value = "'31.35"
$ws.Range("E15").Value = "  -5.40%  "

$ws.Range("D16").Value = "3.475.53"
$ws.Range("E16").Value = "  -2.54%  "

$ws.Range("D17").Value = "66.960.88"
$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("E18").Value = "  +0.22%  "

$ws.Range("D19").Value = "'6.43"
$ws.Range("E19").Value = "  -4.82%  "

$ws.Range("D20").Value = "'15.30"
$ws.Range("E20").Value = "  -4.35%  "

$ws.Range("D21").Value = "'10.05"
$ws.Range("E21").Value = "  +0.62%  "

$ws.Range("D22").Value = "'433.75"
$ws.Range("E22").Value = "  -4.57%  "

$ws.Range("D23").Value = "'0.606"
$ws.Range("E23").Value = "  -5.76%  "

$ws.Range("D24").Value = "'78.97"
$ws.Range("E24").Value = "  +0.44%  "

$ws.Range("E25").Value = "  +0.15%  "

$ws.Range("D26").Value = "3.616.83"
$ws.Range("E26").Value = "  -2.48%  "

$ws.Range("D27").Value = "'0.0000121"
$ws.Range("E27").Value = "  -7.75%  "

$ws.Range("D28").Value = "'9.81"
$ws.Range("E28").Value = "  -6.96%  "

$ws.Range("D29").Value = "'8.39"
$ws.Range("E29").Value = "  -8.42%  "

$ws.Range("D31").Value = "'1.59"
$ws.Range("E31").Value = "  -6.70%  "

$ws.Range("B32").Value = "Kaspa"
$ws.Range("C32").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D32").Value = "'0.167"
$ws.Range("E32").Value = "  -1.99%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").Value = "'1.00"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "'25.31"
$ws.Range("E34").Value = "  -3.14%  "

$ws.Range("D35").Value = "3.470.66"
$ws.Range("E35").Value = "  -2.44%  "

$ws.Range("D36").Value = "'5.93"
$ws.Range("E36").Value = "  -6.47%  "

$ws.Range("D37").Value = "'1.80"
$ws.Range("E37").Value = "  -6.29%  "

$ws.Range("E38").Value = "  +0.02%  "

$ws.Range("D39").Value = "'7.90"
$ws.Range("E39").Value = "  -4.45%  "

$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.15%  "

$ws.Range("D41").Value = "'173.96"
$ws.Range("E41").Value = "  -3.06%  "

$ws.Range("D42").Value = "'0.0882"
$ws.Range("E42").Value = "  -4.00%  "

$ws.Range("D43").Value = "'2.09"
$ws.Range("E43").Value = "  -12.41%  "

$ws.Range("D44").Value = "'5.40"
$ws.Range("E44").Value = "  -3.45%  "

$ws.Range("E45").Value = "  -0.44%  "

$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").Value = "'46.36"
$ws.Range("E46").Value = "  +0.27%  "

$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'28.93"
$ws.Range("E47").Value = "  -6.72%  "

$ws.Range("E48").Value = "  -7.79%  "

$ws.Range("D49").Value = "'7.45"
$ws.Range("E49").Value = "  -4.46%  "

$ws.Range("D50").Value = "'2.43"
$ws.Range("E50").Value = "  -8.77%  "

$ws.Range("D51").Value = "'0.978"
$ws.Range("E51").Value = "  -4.12%  "

